$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 5 (Scenario 4 - Update place) ---
# Clone row 4's formatting down into row 5 first so the new cells pick up
# the same styles used by the other data rows (thin-box border for A:G,
# plain for H, boxed border for I:Q).
$ws.Range("A4:Q4").Copy()
$ws.Range("A5:Q5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A5").Value = "Scenario 4"

# --- Fix up the descriptions on existing rows 3 and 4 ---
# (Row 3 was a copy/paste of the "add" row -> now the Get-place scenario;
#  Row 4 was a copy/paste of the "add" row -> now the Delete-place scenario.)
$ws.Range("B3").Value = "Get place API call"
$ws.Range("B4").Value = "Delete place API call"

# --- Continue filling in the rest of row 5 ---
$ws.Range("B5").Value = "Update place API call"
$ws.Range("C5").Value = "application/json"
$ws.Range("D5").Value = "/maps/api/place/update/json"
$ws.Range("E5").Value = "https://rahulshettyacademy.com"
$ws.Range("F5").Value = "rahulshettyacademy.com"
$ws.Range("G5").Value = "qaclick123"
$ws.Range("H5").Value = -38.383493999999999
$ws.Range("I5").Value = 33.427362000000002
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = "Frontline House"
$ws.Range("L5").Value = "(+91) 929 875 0767"
$ws.Range("M5").Value = $ws.Range("M4").Text
$ws.Range("N5").Value = "Shoe park_shop"
$ws.Range("O5").Value = "http://google.com"
$ws.Range("P5").Value = "English"
$ws.Range("Q5").Value = "ef6ed47dffcf1a24b70ef776662f2bc1"

# --- Hyperlink on the new row's origin cell ---
$ws.Hyperlinks.Add($ws.Range("E5"), "https://rahulshettyacademy.com/")
# Adding the hyperlink re-styles E5 with the built-in Hyperlink style;
# restore the plain bordered look used by the rest of column E.
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the saved view/selection on the sheet ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L5").Select()
